# Re-ran the averaged-intensities notebook, this time including the new
# spiral sampling schemes (Spiral-90deg-10rot-5space, Spiral-90deg-15rot-5space,
# Spiral-90deg-10rot-3space). The scheme table grows from 14 to 17 rows and
# Gaussian-Quadrature moves up next to the 3 new Spiral rows; all downstream
# rows are rewritten with freshly recomputed averaged-intensity values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/bold/centered "id" style (used by column A, e.g. A16)
# down onto the 3 new rows (17-19) before writing values into them.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Full A3:M19 data block (17 schemes x 13 columns: id, name, 11 HKL/pair
# intensity columns), rebuilt in the new row order.
$arr = New-Object 'object[,]' 17,13
$arr[0,0] = 1
$arr[0,1] = "ND Single"
$arr[0,2] = 0.9399999999999999
$arr[0,3] = 1.08
$arr[0,4] = 0.99
$arr[0,5] = 0.9399999999999999
$arr[0,6] = 1.02
$arr[0,7] = 1
$arr[0,8] = 0.98
$arr[0,9] = 1.08
$arr[0,10] = 1.035
$arr[0,11] = 0.9874999999999999
$arr[0,12] = 1.001666666666667
$arr[1,0] = 2
$arr[1,1] = "RD Single"
$arr[1,2] = 1.06
$arr[1,3] = 0.65
$arr[1,4] = 1.09
$arr[1,5] = 1.06
$arr[1,6] = 0.79
$arr[1,7] = 1.26
$arr[1,8] = 1.09
$arr[1,9] = 0.65
$arr[1,10] = 0.8700000000000001
$arr[1,11] = 0.965
$arr[1,12] = 0.9899999999999999
$arr[2,0] = 3
$arr[2,1] = "TD Single"
$arr[2,2] = 1.11
$arr[2,3] = 0.75
$arr[2,4] = 1.02
$arr[2,5] = 1.11
$arr[2,6] = 0.87
$arr[2,7] = 1.06
$arr[2,8] = 1.05
$arr[2,9] = 0.75
$arr[2,10] = 0.885
$arr[2,11] = 0.9975000000000001
$arr[2,12] = 0.9766666666666667
$arr[3,0] = 4
$arr[3,1] = "Morris"
$arr[3,2] = 1.01
$arr[3,3] = 0.77
$arr[3,4] = 1.06
$arr[3,5] = 1.01
$arr[3,6] = 0.87
$arr[3,7] = 1.18
$arr[3,8] = 1.05
$arr[3,9] = 0.77
$arr[3,10] = 0.915
$arr[3,11] = 0.9624999999999999
$arr[3,12] = 0.9899999999999999
$arr[4,0] = 5
$arr[4,1] = "Ring Perpendicular to ND"
$arr[4,2] = 1.002739726027397
$arr[4,3] = 1.001369863013699
$arr[4,4] = 0.9878082191780821
$arr[4,5] = 1.002739726027397
$arr[4,6] = 1.003698630136986
$arr[4,7] = 0.9717808219178082
$arr[4,8] = 0.9913698630136987
$arr[4,9] = 1.001369863013699
$arr[4,10] = 0.9945890410958904
$arr[4,11] = 0.9986643835616439
$arr[4,12] = 0.9931278538812786
$arr[5,0] = 6
$arr[5,1] = "Ring Perpendicular to RD"
$arr[5,2] = 1.032105263157895
$arr[5,3] = 0.8589473684210527
$arr[5,4] = 1.024210526315789
$arr[5,5] = 1.032105263157895
$arr[5,6] = 0.9221052631578948
$arr[5,7] = 1.074736842105263
$arr[5,8] = 1.028421052631579
$arr[5,9] = 0.8589473684210527
$arr[5,10] = 0.9415789473684211
$arr[5,11] = 0.9868421052631579
$arr[5,12] = 0.9900877192982457
$arr[6,0] = 7
$arr[6,1] = "Ring Perpendicular to TD"
$arr[6,2] = 0.9905263157894737
$arr[6,3] = 0.93
$arr[6,4] = 1.019473684210526
$arr[6,5] = 0.9905263157894737
$arr[6,6] = 0.9494736842105264
$arr[6,7] = 1.071052631578947
$arr[6,8] = 1.014210526315789
$arr[6,9] = 0.93
$arr[6,10] = 0.9747368421052631
$arr[6,11] = 0.9826315789473685
$arr[6,12] = 0.9957894736842107
$arr[7,0] = 8
$arr[7,1] = "Gaussian-Quadrature"
$arr[7,2] = 0.9971534375439213
$arr[7,3] = 1.00908493527576
$arr[7,4] = 0.9901014836435641
$arr[7,5] = 0.9971534375439213
$arr[7,6] = 1.001074986944671
$arr[7,7] = 0.9811833402172897
$arr[7,8] = 0.9918300717683775
$arr[7,9] = 1.00908493527576
$arr[7,10] = 0.9995932094596618
$arr[7,11] = 0.9983733235017915
$arr[7,12] = 0.9950713758989305
$arr[8,0] = 9
$arr[8,1] = "Spiral-90deg-10rot-5space"
$arr[8,2] = 0.9927223366655168
$arr[8,3] = 0.9249083018290688
$arr[8,4] = 1.020235861406623
$arr[8,5] = 0.9927223366655168
$arr[8,6] = 0.9474159173647825
$arr[8,7] = 1.072736656152612
$arr[8,8] = 1.01522321488418
$arr[8,9] = 0.9249083018290688
$arr[8,10] = 0.972572081617846
$arr[8,11] = 0.9826472091416815
$arr[8,12] = 0.9955403813837972
$arr[9,0] = 10
$arr[9,1] = "Spiral-90deg-15rot-5space"
$arr[9,2] = 0.9925488203722663
$arr[9,3] = 0.9256919867919757
$arr[9,4] = 1.020038855061651
$arr[9,5] = 0.9925488203722663
$arr[9,6] = 0.9478694004288595
$arr[9,7] = 1.072187902667471
$arr[9,8] = 1.015015511375025
$arr[9,9] = 0.9256919867919757
$arr[9,10] = 0.9728654209268133
$arr[9,11] = 0.9827071206495398
$arr[9,12] = 0.9955587461162082
$arr[10,0] = 11
$arr[10,1] = "Spiral-90deg-10rot-3space"
$arr[10,2] = 0.9926926996768924
$arr[10,3] = 0.9251199895599582
$arr[10,4] = 1.020184540620353
$arr[10,5] = 0.9926926996768924
$arr[10,6] = 0.9475417974180247
$arr[10,7] = 1.072575211627971
$arr[10,8] = 1.015171191936765
$arr[10,9] = 0.9251199895599582
$arr[10,10] = 0.9726522650901558
$arr[10,11] = 0.9826724823835241
$arr[10,12] = 0.995547571806661
$arr[11,0] = 12
$arr[11,1] = "NoRotation-tilt60deg"
$arr[11,2] = 0.9514679999999993
$arr[11,3] = 1.058696000000002
$arr[11,4] = 0.991680000000001
$arr[11,5] = 0.9514679999999993
$arr[11,6] = 1.014340000000001
$arr[11,7] = 1.000455999999999
$arr[11,8] = 0.9839719999999988
$arr[11,9] = 1.058696000000002
$arr[11,10] = 1.025188000000002
$arr[11,11] = 0.9883280000000004
$arr[11,12] = 1.000102
$arr[12,0] = 13
$arr[12,1] = "Rotation-NoTilt"
$arr[12,2] = 0.9437374999999978
$arr[12,3] = 1.08
$arr[12,4] = 0.99
$arr[12,5] = 0.9437374999999978
$arr[12,6] = 1.02
$arr[12,7] = 1
$arr[12,8] = 0.98
$arr[12,9] = 1.08
$arr[12,10] = 1.035
$arr[12,11] = 0.9893687499999989
$arr[12,12] = 1.002289583333333
$arr[13,0] = 14
$arr[13,1] = "Rotation-60detTilt"
$arr[13,2] = 0.9655193791488024
$arr[13,3] = 1.044194979123195
$arr[13,4] = 0.9918701789184027
$arr[13,5] = 0.9655193791488024
$arr[13,6] = 1.0102409080832
$arr[13,7] = 0.9970730475519978
$arr[13,8] = 0.9862887718912008
$arr[13,9] = 1.044194979123195
$arr[13,10] = 1.018032579020799
$arr[13,11] = 0.9917759790848006
$arr[13,12] = 0.9991978774527998
$arr[14,0] = 15
$arr[14,1] = "HexGrid-90degTilt5degRes"
$arr[14,2] = 0.9946215219084183
$arr[14,3] = 0.9947559078995181
$arr[14,4] = 0.994695668562918
$arr[14,5] = 0.9946215219084183
$arr[14,6] = 0.9951444881825868
$arr[14,7] = 0.9945641394354577
$arr[14,8] = 0.9947039883334201
$arr[14,9] = 0.9947559078995181
$arr[14,10] = 0.994725788231218
$arr[14,11] = 0.9946736550698182
$arr[14,12] = 0.9947476190537198
$arr[15,0] = 16
$arr[15,1] = "HexGrid-90degTilt22p5degRes"
$arr[15,2] = 0.997130688889636
$arr[15,3] = 0.991010012114939
$arr[15,4] = 0.9941941151546674
$arr[15,5] = 0.997130688889636
$arr[15,6] = 0.9937293683867972
$arr[15,7] = 0.9942718202339605
$arr[15,8] = 0.9948219014553198
$arr[15,9] = 0.991010012114939
$arr[15,10] = 0.9926020636348032
$arr[15,11] = 0.9948663762622195
$arr[15,12] = 0.9941929843725533
$arr[16,0] = 17
$arr[16,1] = "HexGrid-60degTilt5degRes"
$arr[16,2] = 0.9966918028239566
$arr[16,3] = 0.9886969980975101
$arr[16,4] = 0.995775904112116
$arr[16,5] = 0.9966918028239566
$arr[16,6] = 0.9918368383494247
$arr[16,7] = 0.997957694070406
$arr[16,8] = 0.9963872356812298
$arr[16,9] = 0.9886969980975101
$arr[16,10] = 0.9922364511048131
$arr[16,11] = 0.9944641269643848
$arr[16,12] = 0.9945577455224406
$ws.Range("A3:M19").Value = $arr
